# PreqKpiAgent API workbook update
# - add new "getCaseList" / "getMonRsltList" / "updMonRsltList" endpoints to
#   the "Op" sheet's request-building table
# - add a new scratch sheet "工作表1" listing the fully built URLs

$wb = $excel.ActiveWorkbook
$op = $wb.Worksheets.Item("Op")

# --- new sheet "工作表1" placed after "Op", pre-populated with the three
#     fully-resolved URLs (this mirrors the order the strings were typed in
#     the original edit: scratch sheet first, then the template rows) -----
$new = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $op)
$new.Name = "工作表1"
$new.Range("B3").Value = "http://localhost:8081/updMonRsltList?fabName=FAB18"
$new.Range("B6").Value = "http://localhost:8081/getCaseList"
$new.Range("B8").Value = "http://localhost:8081/getMonRsltList"
$new.Range("B3:B8").Select() | Out-Null

# --- back on "Op": three new request rows --------------------------------
$op.Range("B13").Value = "http://"
$op.Range("C13").Value = "localhost:8081"
$op.Range("F13").Value = "fabName"
$op.Range("G13").Value = "FAB18"
$op.Range("E13").Value = "/updMonRsltList"

$op.Range("B14").Value = "http://"
$op.Range("C14").Value = "localhost:8081"
$op.Range("E14").Value = "/getCaseList"

$op.Range("B15").Value = "http://"
$op.Range("C15").Value = "localhost:8081"
$op.Range("E15").Value = "/getMonRsltList"

# --- extend the URL-building formula in column I down through row 15, and
#     make it also append "?key=value" only when F is not blank ----------
for ($r = 3; $r -le 15; $r++) {
    $op.Range("I$r").Formula = '=CONCATENATE(B' + $r + ',C' + $r + ',D' + $r + ',E' + $r + ', IF(ISBLANK(F' + $r + '), "", CONCATENATE("?",F' + $r + ',"=",G' + $r + ')))'
}

# --- cosmetics: widen column I (auto-fit no longer applies) and move the
#     active selection --------------------------------------------------
$op.Columns("I").ColumnWidth = 38.37
$op.Range("J18").Select() | Out-Null

$op.Select() | Out-Null
